$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Step 1: make room for 1 new row by shifting rows 12..15 down to 13..16.
# Work from bottom to top so we don't overwrite data we still need to read.
# The destination row beyond the original used range (16) must first get a
# style via a format-only paste from a same-shaped source, then the value is
# assigned afterwards - otherwise the new cell ends up with no style index.

$ws.Range("A15:B15").Copy()
$ws.Range("A16:B16").PasteSpecial(-4122)
$ws.Cells.Item(16,1).Value = $ws.Cells.Item(15,1).Value()
$ws.Cells.Item(16,2).Value = $ws.Cells.Item(15,2).Value()

$ws.Cells.Item(15,1).Value = $ws.Cells.Item(14,1).Value()
$ws.Cells.Item(15,2).Value = $ws.Cells.Item(14,2).Value()

$ws.Cells.Item(14,1).Value = $ws.Cells.Item(13,1).Value()
$ws.Cells.Item(14,2).Value = $ws.Cells.Item(13,2).Value()

$ws.Cells.Item(13,1).Value = $ws.Cells.Item(12,1).Value()
$ws.Cells.Item(13,2).Value = $ws.Cells.Item(12,2).Value()

# Step 2: update the metadata values in place (rows 1-10 unchanged except these)
$ws.Cells.Item(3,2).Value = "0.1.7"
$ws.Cells.Item(6,2).Value = "draft"
$ws.Cells.Item(8,2).Value = "2024-11-22T12:33:30-06:00"
$ws.Cells.Item(10,2).Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"

# Step 3: row 11 now becomes the second Contact line, row 12 becomes Jurisdiction
$ws.Cells.Item(11,1).Value = "Contact"
$ws.Cells.Item(11,2).Value = "Bob Milius (bmilius@nmdp.org)"
$ws.Cells.Item(12,1).Value = "Jurisdiction"
$ws.Cells.Item(12,2).Value = ""

Write-Host "done"
